$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing H13 / H14 values
$ws.Range("H13").Value = 19730
$ws.Range("H14").Value = 19881

# Fill in newly populated rows 15-17 (A:I)
$data = @(
    @(14, 0, 6227, 19547, 0, 14601, 0, 20973, 22230),
    @(15, 0, 6249, 19436, 0, 14422, 0, 20862, 22243),
    @(16, 0, 5989, 18881, 0, 14102, 0, 19536, 21365)
)

$startRow = 15
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
